# Add two new columns: I ("I0") and J ("IF") with data for rows 2-37.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (style matches the other header cells, e.g. H1's bold/border/centered style)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-37
$data = @(
    @(2, 9, 9),
    @(3, 9, 9),
    @(4, 7, 7),
    @(5, 8, 8),
    @(6, 7, 7),
    @(7, 7, 7),
    @(8, 7, 7),
    @(9, 8, 8),
    @(10, 4, 5),
    @(11, 7, 7),
    @(12, 10, 11),
    @(13, 8, 8),
    @(14, 8, 8),
    @(15, 7, 8),
    @(16, 6, 6),
    @(17, 7, 7),
    @(18, 7, 8),
    @(19, 7, 7),
    @(20, 7, 7),
    @(21, 8, 8),
    @(22, 7, 7),
    @(23, 3, 4),
    @(24, 10, 10),
    @(25, 6, 6),
    @(26, 7, 7),
    @(27, 7, 7),
    @(28, 8, 8),
    @(29, 8, 8),
    @(30, 8, 8),
    @(31, 6, 6),
    @(32, 9, 9),
    @(33, 7, 7),
    @(34, 8, 8),
    @(35, 4, 4),
    @(36, 6, 6),
    @(37, 7, 7)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
